# Append the new daily row (2025-11-14) to the "Daily 100 Error Counts" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row 28: Date, Total Count, Session Timeout Errors, Errors Requiring Analysis
$ws.Range("A28").Value = 45975
$ws.Range("A28").NumberFormat = "d-mmm-yy"
$ws.Range("B28").Value = 614
$ws.Range("C28").Value = 28
$ws.Range("D28").Value = 586

# Mirror the author's workflow: select the freshly-entered row before saving,
# which is what moves the sheet's recorded selection to A28:D28.
$ws.Range("A28:D28").Select() | Out-Null
